$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 6240.749999999998
$ws.Range("C2").Value = 1593.382978723404
$ws.Range("D2").Value = 3.900468749999999
$ws.Range("E2").Value = 3.0003605769230757

$ws.Range("B3").Value = 5991.119999999997
$ws.Range("C3").Value = 1529.6476595744675
$ws.Range("D3").Value = 3.7444499999999983
$ws.Range("E3").Value = 2.880346153846153

$ws.Range("B4").Value = 482.88888888888886
$ws.Range("C4").Value = 123.29078014184398
$ws.Range("D4").Value = 0.30180555555555555
$ws.Range("E4").Value = 0.23215811965811964

$ws.Range("B5").Value = 12714.758888888884
$ws.Range("C5").Value = 3246.3214184397157
$ws.Range("D5").Value = 7.946724305555552
$ws.Range("E5").Value = 6.112864850427349

$ws.Range("B7").Value = 2819.9999999999995
$ws.Range("D7").Value = 1.7624999999999997
$ws.Range("E7").Value = 1.3557692307692306

$ws.Range("B8").Value = 1057.4999999999998
$ws.Range("D8").Value = 0.6609374999999998
$ws.Range("E8").Value = 0.5084134615384615

$ws.Range("B9").Value = 3877.4999999999995
$ws.Range("D9").Value = 2.4234375
$ws.Range("E9").Value = 1.864182692307692

$ws.Range("B11").Value = 3684.8341009758096
$ws.Range("C11").Value = 940.8087066321217
$ws.Range("D11").Value = 2.303021313109881
$ws.Range("E11").Value = 1.77155485623837

$ws.Range("C13").Value = 129.81125106382976
$ws.Range("D13").Value = 0.3177671249999999
$ws.Range("E13").Value = 0.2444362499999999

$ws.Range("B14").Value = 1691.6602304037292
$ws.Range("C14").Value = 431.9132503158458

$ws.Range("B15").Value = 1403.9999999999995
$ws.Range("C15").Value = 358.46808510638294
$ws.Range("D15").Value = 0.8774999999999997
$ws.Range("E15").Value = 0.6749999999999998

$ws.Range("B16").Value = 3604.0876304037283
$ws.Range("C16").Value = 920.1925864860584
$ws.Range("D16").Value = 2.25255476900233
$ws.Range("E16").Value = 1.7327344376941

$ws.Range("B18").Value = 1877.6660181615362
$ws.Range("C18").Value = 479.404089743371
$ws.Range("D18").Value = 1.1735412613509602
$ws.Range("E18").Value = 0.9027240471930462

$ws.Range("B19").Value = 2914.245316464084
$ws.Range("C19").Value = 744.0626339908301
$ws.Range("D19").Value = 1.8214033227900526
$ws.Range("E19").Value = 1.4010794790692713

$ws.Range("B20").Value = 4901.598702989739
$ws.Range("C20").Value = 1251.472009273976
$ws.Range("D20").Value = 3.063499189368587
$ws.Range("E20").Value = 2.3565378379758357

$ws.Range("B23").Value = 28782.77932325816
$ws.Range("C23").Value = 7348.794720831872
$ws.Range("D23").Value = 17.98923707703635
$ws.Range("E23").Value = 13.837874674643347

$ws.Range("B25").Value = 16068.020434369277
$ws.Range("C25").Value = 4102.473302392156
$ws.Range("D25").Value = 10.042512771480798
$ws.Range("E25").Value = 7.725009824215999
